# Update the workbook per the author's edit:
#  - progress (%) figures in column I for rows 9-16 were revised
#  - the sheet view had scrolled down and the active selection moved

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("개발목록")
$ws.Activate()

# Revised progress percentages (stored as fractions, displayed as %)
$ws.Range("I9").Value  = 0.8
$ws.Range("I10").Value = 0.5
$ws.Range("I11").Value = 0.4
$ws.Range("I13").Value = 0.4
$ws.Range("I14").Value = 0.4
$ws.Range("I15").Value = 0.4
$ws.Range("I16").Value = 0.4

# View state: the frozen header (top 2 rows) stays in place, but the
# window had been scrolled down so row 6 is the first visible row below
# the freeze, and the active cell/selection moved to I17.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I17").Select()
